$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @{
    2 = "Boitumelo"
    3 = "Frank"
    4 = "None"
    5 = "An"
    6 = "Elsa"
    7 = "None"
    8 = "Kevin P"
    9 = "Miro"
    10 = "Patrick"
    11 = "Edoardo"
    12 = "None"
    13 = "Beatriz"
    14 = "Fatemeh"
    15 = "Karthika"
    16 = "Nicole"
    17 = "Maxim"
    18 = "Miriam"
    19 = "Andrii"
    20 = "Manel"
    21 = "Jean"
    22 = "Celina"
    23 = "Patrycja"
    24 = "Aleksander"
    25 = "Jessica"
    26 = "Dhanya"
    27 = "Oscar"
    28 = "Yassine"
    29 = "Mohamad"
    30 = "Nina"
    31 = "David"
    32 = "Vera"
    33 = "Olha"
    34 = "Kevin J"
    35 = "Therese"
    36 = "Stef"
    37 = "Imad"
}

foreach ($row in $names.Keys) {
    $ws.Range("B$row").Value = $names[$row]
}
